# Update dSF (column F) values to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F4").Value = -5
$ws.Range("F9").Value = -10
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = -3
$ws.Range("F17").Value = 3
$ws.Range("F21").Value = 3
